# Handback status report regeneration.
#
# Running the handback report again refreshes the "Correspond Handoff
# Datetime" (col D) and "Correspond Handback DateTime" (col G) timestamps
# for the most recently processed file (row 2) on each language sheet,
# while the already-handled row (row 3) keeps its historical timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-22 03:25:11"
$wsZhCn.Range("G2").Value = "2016-02-22 03:25:58"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-22 03:25:24"
$wsDeDe.Range("G2").Value = "2016-02-22 03:26:17"
